$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width changes (COM ColumnWidth is stored internally with a fixed
# offset vs. the OOXML <col width> value, so subtract that offset up front)
$ws.Columns.Item(3).ColumnWidth = 1.3072916666666665
$ws.Columns.Item(4).ColumnWidth = 1.3072916666666665
$ws.Columns.Item(15).ColumnWidth = 3.8776041666666665
$ws.Columns.Item(16).ColumnWidth = 3.8776041666666665
$ws.Columns.Item(17).ColumnWidth = 3.8776041666666665

# Cell value changes in row 1
$ws.Range("C1").Value = 8
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 25
$ws.Range("F1").Value = 6
$ws.Range("G1").Value = 16
$ws.Range("H1").Value = 22
$ws.Range("I1").Value = 16
$ws.Range("J1").Value = 19
$ws.Range("K1").Value = 26
$ws.Range("L1").Value = 13
$ws.Range("M1").Value = 0.040999999999999995
$ws.Range("N1").Value = 0.033000000000000002
$ws.Range("O1").Value = 0.059999999999999998
$ws.Range("P1").Value = 0.089999999999999997
$ws.Range("Q1").Value = 0.069999999999999993
